$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

$urlMd1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f0fa35286659f871b2fbc65d03c3a99bcd34508b/e2e/0a496b9b-ebc7-42b6-836e-119fc7c40731.md"
$urlMd2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f0fa35286659f871b2fbc65d03c3a99bcd34508b/e2e/e0f8d963-78c4-442e-98a7-99f85ac027a5.md"

# --- Status text: "In Translation" -> "Handed back: in sync with en-US" everywhere it shows up ---
$newStatus = "Handed back: in sync with en-US"

$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Overview: bump the "Latest HO Xliff Generate Date" to reflect the handback run ---
# (G2/G3 already show 2016-09-04 22:26:59, unchanged by the diff)

# --- zh-cn sheet: data updates for rows 2/3 ---
$wsZhCn.Range("J2").Value = "0a496b9b-ebc7-42b6-836e-119fc7c40731.8d1f6f3c2b638b8aa3ed7c26f3926eb94bbadf1d.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-09-04 22:27:56"
$wsZhCn.Range("J3").Value = "e0f8d963-78c4-442e-98a7-99f85ac027a5.c804adcea41ad230967ab2cd309488c58c65bf7f.zh-cn.xlf"
$wsZhCn.Range("K3").Value = "2016-09-04 22:27:56"

# --- de-de sheet: data updates for rows 2/3 ---
$wsDeDe.Range("J2").Value = "0a496b9b-ebc7-42b6-836e-119fc7c40731.8d1f6f3c2b638b8aa3ed7c26f3926eb94bbadf1d.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-09-04 22:28:09"
$wsDeDe.Range("J3").Value = "e0f8d963-78c4-442e-98a7-99f85ac027a5.c804adcea41ad230967ab2cd309488c58c65bf7f.de-de.xlf"
$wsDeDe.Range("K3").Value = "2016-09-04 22:28:09"

# --- Rebuild the hyperlink list on each sheet so the new "Latest Target File"
#     links (I2/I3) interleave with the existing "Source File Name" links
#     (A2/A3) in A2, I2, A3, I3 order, matching a regenerated report. ---
$wsZhCn.Cells.Hyperlinks.Delete()
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A2"), $urlMd1, "", "", "0a496b9b-ebc7-42b6-836e-119fc7c40731.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $urlMd1, "", "", "0a496b9b-ebc7-42b6-836e-119fc7c40731.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $urlMd2, "", "", "e0f8d963-78c4-442e-98a7-99f85ac027a5.md")
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I3"), $urlMd2, "", "", "e0f8d963-78c4-442e-98a7-99f85ac027a5.md")

$wsDeDe.Cells.Hyperlinks.Delete()
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A2"), $urlMd1, "", "", "0a496b9b-ebc7-42b6-836e-119fc7c40731.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $urlMd1, "", "", "0a496b9b-ebc7-42b6-836e-119fc7c40731.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $urlMd2, "", "", "e0f8d963-78c4-442e-98a7-99f85ac027a5.md")
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I3"), $urlMd2, "", "", "e0f8d963-78c4-442e-98a7-99f85ac027a5.md")

# --- Column width adjustments to fit the newly-populated long text/hyperlink columns ---
$wsOverview.Range("E1").ColumnWidth = 29.9777047293527
$wsOverview.Range("F1").ColumnWidth = 29.9777047293527

$wsZhCn.Range("C1").ColumnWidth = 29.9777047293527
$wsZhCn.Range("I1").ColumnWidth = 40
$wsZhCn.Range("J1").ColumnWidth = 40

$wsDeDe.Range("C1").ColumnWidth = 29.9777047293527
$wsDeDe.Range("I1").ColumnWidth = 40
$wsDeDe.Range("J1").ColumnWidth = 40

Write-Host "Handback report generated"
